# The pricing data refresh removed every individual service line item
# (EC2/S3/RDS/ALB) from the cost report and rolled the total down to
# $0.00 since no costed services remain. Rows 1 and 3 (title + header)
# are untouched; row 5 keeps its row number but is overwritten to become
# the new "Total Estimated Monthly Cost" / "$0.00" row; rows 4, 6, 7 and
# 9 are wiped out entirely (not shifted) so the sheet's used range
# shrinks to A1:B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the EC2 line (row 4) completely.
$ws.Rows(4).ClearContents()

# Row 5 (previously "S3" / "$0.02") becomes the new totals row.
# Force the amount cell to text first so "$0.00" is stored as a literal
# string (matching the other cost cells) instead of being auto-detected
# by Excel as a currency number.
$ws.Range("B5").NumberFormat = "@"
$ws.Range("A5").Value = "Total Estimated Monthly Cost"
$ws.Range("B5").Value = "$0.00"

# Remove the remaining RDS, ALB and old Total rows completely.
$ws.Rows(6).ClearContents()
$ws.Rows(7).ClearContents()
$ws.Rows(9).ClearContents()
